# Generate Report for Archive
#
# The localization status changed from "Ready for handoff" to
# "In Translation" for the tracked file. Update every sheet that shows
# this status (Overview summary columns + each per-locale detail sheet),
# then shrink the now-narrower status columns to match the shorter text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: per-locale status columns (E = zh-cn, F = de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn detail sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de detail sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
